$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("M18").ClearContents()
$ws.Range("H32").Value = 989.1818
$ws.Range("I32").Value = 900
$ws.Range("J32").Value = 998.1
$ws.Range("K32").Value = 900
$ws.Range("L32").Value = 998.1
$ws.Range("M32").Value = -574
$ws.Range("N32").Value = -1650.1
$ws.Range("H62").Value = 6878
$ws.Range("I62").Value = 6254.4287
$ws.Range("K62").Value = 6254.4287
$ws.Range("M62").Value = -5630.4287
$ws.Range("H65").Value = 6878
$ws.Range("I65").Value = 6254.4287
$ws.Range("K65").Value = 31272.1435
$ws.Range("M65").Value = -28152.1435
$ws.Range("H69").Value = 2500
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()
$ws.Range("H72").Value = 2500
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()
$ws.Range("H107").Value = 879.2727
$ws.Range("I107").Value = 1075
$ws.Range("J107").Value = 767.4286
$ws.Range("K107").Value = 1075
$ws.Range("L107").Value = 767.4286
$ws.Range("M107").Value = 845
$ws.Range("N107").Value = -4607.4286
$ws.Range("H113").Value = 7112.5713
$ws.Range("I113").Value = 5322.125
$ws.Range("J113").Value = 9499.833000000001
$ws.Range("K113").Value = 5322.125
$ws.Range("L113").Value = 9499.833000000001
$ws.Range("M113").Value = -2068.125
$ws.Range("N113").Value = -16007.833
$ws.Range("H116").Value = 3905.818
$ws.Range("I116").Value = 2537.4
$ws.Range("K116").Value = 2537.4
$ws.Range("M116").Value = 904.5999999999999
$ws.Range("H132").Value = 3096.5454
$ws.Range("I132").Value = 3981.375
$ws.Range("K132").Value = 11944.125
$ws.Range("M132").Value = -9414.125
$ws.Range("H137").Value = 1873
$ws.Range("I137").Value = 1355.5
$ws.Range("K137").Value = 4066.5
$ws.Range("M137").Value = -1516.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 2898
$ws.Range("I102").Value = 2872.5
$ws.Range("K102").Value = 2872.5
$ws.Range("M102").Value = -1250.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1166.625
$ws.Range("I20").Value = 892
$ws.Range("J20").Value = 1331.4
$ws.Range("K20").Value = 892
$ws.Range("L20").Value = 1331.4
$ws.Range("M20").Value = -645
$ws.Range("N20").Value = -1825.4
$ws.Range("H22").Value = 114
$ws.Range("I22").Value = 90
$ws.Range("K22").Value = 90
$ws.Range("M22").Value = 83
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").ClearContents()
$ws.Range("N95").Value = 0
$ws.Range("H100").Value = 29573.834
$ws.Range("J100").Value = 29573.834
$ws.Range("L100").Value = 29573.834
$ws.Range("N100").Value = -31737.834
$ws.Range("H106").Value = 8589.75
$ws.Range("J106").Value = 8589.75
$ws.Range("L106").Value = 8589.75
$ws.Range("N106").Value = -11113.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 10536.4
$ws.Range("H31").Value = 3242
$ws.Range("I31").Value = 2940.2
$ws.Range("K31").Value = 2940.2
$ws.Range("M31").Value = -2645.2
$ws.Range("H34").Value = 3242
$ws.Range("I34").Value = 2940.2
$ws.Range("K34").Value = 2940.2
$ws.Range("M34").Value = -2738.2
$ws.Range("H38").Value = 7092
$ws.Range("I38").Value = 8538
$ws.Range("J38").Value = 4200
$ws.Range("K38").Value = 8538
$ws.Range("L38").Value = 4200
$ws.Range("M38").Value = -8161
$ws.Range("N38").Value = -4954
$ws.Range("H43").Value = 15000
$ws.Range("J43").Value = 15000
$ws.Range("L43").Value = 15000
$ws.Range("N43").Value = -15368
$ws.Range("H46").Value = 7092
$ws.Range("I46").Value = 8538
$ws.Range("J46").Value = 4200
$ws.Range("K46").Value = 8538
$ws.Range("L46").Value = 4200
$ws.Range("M46").Value = -8327
$ws.Range("N46").Value = -4622
$ws.Range("H48").Value = 37589
$ws.Range("J48").Value = 79575.25
$ws.Range("L48").Value = 79575.25
$ws.Range("N48").Value = -80527.25
$ws.Range("H58").Value = 3590.0908
$ws.Range("I58").Value = 2262
$ws.Range("K58").Value = 2262
$ws.Range("M58").Value = -2059
$ws.Range("H88").Value = 29740.5
$ws.Range("J88").Value = 29740.5
$ws.Range("L88").Value = 29740.5
$ws.Range("N88").Value = -30552.5
$ws.Range("H91").Value = 29740.5
$ws.Range("J91").Value = 29740.5
$ws.Range("L91").Value = 29740.5
$ws.Range("N91").Value = -32548.5
$ws.Range("H101").Value = 15000
$ws.Range("J101").Value = 15000
$ws.Range("L101").Value = 15000
$ws.Range("N101").Value = -21490
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").ClearContents()
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = 0
$ws.Range("H136").Value = 3590.0908
$ws.Range("I136").Value = 2262
$ws.Range("K136").Value = 6786
$ws.Range("M136").Value = -4236

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 3497.5
$ws.Range("I8").Value = 3497.5
$ws.Range("K8").Value = 10492.5
$ws.Range("M8").Value = -10353.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 670000
$ws.Range("J7").Value = 510000
$ws.Range("L7").Value = 510000
$ws.Range("N7").Value = -510224
$ws.Range("H8").Value = 670000
$ws.Range("J8").Value = 510000
$ws.Range("L8").Value = 510000
$ws.Range("N8").Value = -510278

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 12333
$ws.Range("J20").Value = 18000
$ws.Range("L20").Value = 18000
$ws.Range("N20").Value = -18452
$ws.Range("H46").Value = 4343.625
$ws.Range("I46").Value = 2874.5
$ws.Range("J46").Value = 4833.3335
$ws.Range("K46").Value = 2874.5
$ws.Range("L46").Value = 4833.3335
$ws.Range("M46").Value = -2686.5
$ws.Range("N46").Value = -5209.3335
$ws.Range("H61").Value = 1956.6
$ws.Range("I61").Value = 2145.75
$ws.Range("K61").Value = 2145.75
$ws.Range("M61").Value = -1943.75
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").ClearContents()
$ws.Range("N95").Value = 0
$ws.Range("H100").Value = 2309.0908
$ws.Range("J100").Value = 3200
$ws.Range("L100").Value = 3200
$ws.Range("N100").Value = -4282
$ws.Range("H101").Value = 25362
$ws.Range("J101").Value = 25362
$ws.Range("L101").Value = 25362
$ws.Range("N101").Value = -31852
$ws.Range("H113").Value = 1956.6
$ws.Range("I113").Value = 2145.75
$ws.Range("K113").Value = 2145.75
$ws.Range("M113").Value = 24.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 35635
$ws.Range("J69").Value = 35635
$ws.Range("L69").Value = 35635
$ws.Range("N69").Value = -37133
$ws.Range("H72").Value = 35635
$ws.Range("J72").Value = 35635
$ws.Range("L72").Value = 106905
$ws.Range("N72").Value = -114393
$ws.Range("H96").Value = 1197.25
$ws.Range("I96").Value = 1197.25
$ws.Range("K96").Value = 1197.25
$ws.Range("M96").Value = 175.75
$ws.Range("H107").Value = 1607.9546
$ws.Range("I107").Value = 1586.2142
$ws.Range("J107").Value = 1646
$ws.Range("K107").Value = 4758.642599999999
$ws.Range("L107").Value = 4938
$ws.Range("M107").Value = -2838.642599999999
$ws.Range("N107").Value = -8778
